$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row for "un_franzosa_ControlvsCD_ConvCD" above the existing
# "un_franzosa_ControlvsCD_Fp" row (currently row 9), shifting subsequent
# rows down by one.
$ws.Rows("9:9").Insert()
$ws.Range("A9").Value = "un_franzosa_ControlvsCD_ConvCD"
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0.33
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.67
$ws.Range("H9").Value = 0.67

# Insert a new row for "un_franzosa_ControlvsUC_ConvUC" above the existing
# "un_franzosa_ControlvsUC_Fp" row (now at row 14 after the previous
# insertion), shifting subsequent rows down by one.
$ws.Rows("14:14").Insert()
$ws.Range("A14").Value = "un_franzosa_ControlvsUC_ConvUC"
$ws.Range("B14").Value = 0
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0.67
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.33
$ws.Range("H14").Value = 0.33
